# The authoring edit removed a redundant/duplicate picture ("Imagen 8",
# shape id 9) from slide 6 - it was an untouched-up duplicate sitting behind
# the final "Imagen 9" picture (which has the background removed +
# AI-generated description). Everything else on the slide is left as-is.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Imagen 8" -and $sh.Id -eq 9) {
        $sh.Delete()
    }
}
